$wb = $excel.ActiveWorkbook

# Sheet 1 (土地)
$ws = $wb.Worksheets.Item(1)
$ws.Range('B1').Value = 'name'
$ws.Range('C1').Value = 'area'
$ws.Range('D1').Value = 'share_portion'
$ws.Range('E1').Value = 'owner'
$ws.Range('F1').Value = 'register_date'
$ws.Range('G1').Value = 'register_reason'
$ws.Range('H1').Value = 'acquire_value'
$ws.Range('I1').Value = 'property_category'
$ws.Range('J1').Value = 'category'
$ws.Range('K1').Value = 'date'
$ws.Range('L1').Value = 'legislator_name'
$ws.Range('M1').Value = 'legislator_id'
$ws.Range('N1').Value = 'source_file'
$ws.Range('O1').Value = 'index'
$ws.Range('B2').Value = '新北市土城區永和段05570000地號'
$ws.Range('D2').Value = '全部'
$ws.Range('E2').Value = '盧嘉辰'
$ws.Range('F2').Value = '70年06月10日'
$ws.Range('G2').Value = '買賣'
$ws.Range('H2').Value = '(超過五年）'
$ws.Range('I2').Value = 'land'
$ws.Range('J2').Value = 'normal'
$ws.Range('K2').Value = '2012-04-12'
$ws.Range('L2').Value = '盧嘉辰'
$ws.Range("M2").Value = 1715
$ws.Range('N2').Value = 'tmp79201'
$ws.Range("O2").Value = 13
$ws.Range('B3').Value = '新北市土城區永和段05580000地號'
$ws.Range('D3').Value = '全部'
$ws.Range('E3').Value = '盧嘉辰'
$ws.Range('F3').Value = '73年06月15日'
$ws.Range('G3').Value = '買賣'
$ws.Range('H3').Value = '(超過五年）'
$ws.Range('I3').Value = 'land'
$ws.Range('J3').Value = 'normal'
$ws.Range('K3').Value = '2012-04-12'
$ws.Range('L3').Value = '盧嘉辰'
$ws.Range("M3").Value = 1715
$ws.Range('N3').Value = 'tmp79201'
$ws.Range("O3").Value = 14
$ws.Range('B4').Value = '新北市上城區永和段07010000地號'
$ws.Range('D4').Value = '36分之1'
$ws.Range('E4').Value = '盧嘉辰'
$ws.Range('F4').Value = '98年09月07日'
$ws.Range('G4').Value = '受贈'
$ws.Range('I4').Value = 'land'
$ws.Range('J4').Value = 'normal'
$ws.Range('K4').Value = '2012-04-12'
$ws.Range('L4').Value = '盧嘉辰'
$ws.Range("M4").Value = 1715
$ws.Range('N4').Value = 'tmp79201'
$ws.Range("O4").Value = 15
$ws.Range('B5').Value = '新北市土城區永和段07020000地號'
$ws.Range('D5').Value = '36分之1'
$ws.Range('E5').Value = '盧嘉辰'
$ws.Range('F5').Value = '98年09月07H'
$ws.Range('G5').Value = '受贈'
$ws.Range('I5').Value = 'land'
$ws.Range('J5').Value = 'normal'
$ws.Range('K5').Value = '2012-04-12'
$ws.Range('L5').Value = '盧嘉辰'
$ws.Range("M5").Value = 1715
$ws.Range('N5').Value = 'tmp79201'
$ws.Range("O5").Value = 16
$ws.Range('B6').Value = '新北市土城區永和段07160000地號'
$ws.Range('D6').Value = '18分之1'
$ws.Range('E6').Value = '盧嘉辰'
$ws.Range('F6').Value = '98年09月07日'
$ws.Range('G6').Value = '受贈'
$ws.Range('I6').Value = 'land'
$ws.Range('J6').Value = 'normal'
$ws.Range('K6').Value = '2012-04-12'
$ws.Range('L6').Value = '盧嘉辰'
$ws.Range("M6").Value = 1715
$ws.Range('N6').Value = 'tmp79201'
$ws.Range("O6").Value = 17
$ws.Range('B7').Value = '新北市上城區永和段07320000地號'
$ws.Range('D7').Value = '6分之1'
$ws.Range('E7').Value = '盧嘉辰'
$ws.Range('F7').Value = '98年09月07H'
$ws.Range('G7').Value = '受贈'
$ws.Range('I7').Value = 'land'
$ws.Range('J7').Value = 'normal'
$ws.Range('K7').Value = '2012-04-12'
$ws.Range('L7').Value = '盧嘉辰'
$ws.Range("M7").Value = 1715
$ws.Range('N7').Value = 'tmp79201'
$ws.Range("O7").Value = 18
$ws.Range('B8').Value = '新北市土城區永和段07420000地號'
$ws.Range('D8').Value = '6分之1'
$ws.Range('E8').Value = '盧嘉辰'
$ws.Range('F8').Value = '98年09月07日'
$ws.Range('G8').Value = '受贈'
$ws.Range('I8').Value = 'land'
$ws.Range('J8').Value = 'normal'
$ws.Range('K8').Value = '2012-04-12'
$ws.Range('L8').Value = '盧嘉辰'
$ws.Range("M8").Value = 1715
$ws.Range('N8').Value = 'tmp79201'
$ws.Range("O8").Value = 19
$ws.Range('B9').Value = '新北市土城區建安段00020000地號'
$ws.Range('D9').Value = '6分之1'
$ws.Range('E9').Value = '盧嘉辰'
$ws.Range('F9').Value = '98年09月07日'
$ws.Range('G9').Value = '受贈'
$ws.Range('I9').Value = 'land'
$ws.Range('J9').Value = 'normal'
$ws.Range('K9').Value = '2012-04-12'
$ws.Range('L9').Value = '盧嘉辰'
$ws.Range("M9").Value = 1715
$ws.Range('N9').Value = 'tmp79201'
$ws.Range("O9").Value = 20

# Sheet 2 (建物)
$ws = $wb.Worksheets.Item(2)
$ws.Range('B2').Value = '新北市上城區永和段00006029建號'
$ws.Range('F2').Value = '73年10月23日'
$ws.Range('B3').Value = '新北市土城區永和段00006030建號'
$ws.Range('F3').Value = '73年10月23U'

# Sheet 3 (汽車)
$ws = $wb.Worksheets.Item(3)
$ws.Range('E2').Value = '87年08月11H'
$ws.Range('E3').Value = '97年01月27R'
$ws.Range('E4').Value = '101年02月13曰'

# Sheet 5 (債務)
$ws = $wb.Worksheets.Item(5)
$ws.Range('D2').Value = '永豐銀行臺北市中正區南昌路'
$ws.Range('F2').Value = '91年07月17日'

# Sheet 6 (事業投資)
$ws = $wb.Worksheets.Item(6)
$ws.Range('F2').Value = '68年02月07日'
$ws.Range('F3').Value = '68年02月07日'
